$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.523.63"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "3.100.24"
$ws.Range("E3").Value = "  -2.03%  "
$ws.Range("D4").Value = "'0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").Value = "'230.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.44%  "
$ws.Range("D6").Value = "'625.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("D8").Value = "'0.363"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.73%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "3.100.23"
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("D11").Value = "'0.725"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.01%  "
$ws.Range("D12").Value = "'0.196"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").Value = "'36.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.01%  "
$ws.Range("E14").Value = "  -1.62%  "
$ws.Range("E15").Value = "  -1.96%  "
$ws.Range("D16").Value = "90.592.13"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "3.693.35"
$ws.Range("E17").Value = "  -1.80%  "
$ws.Range("D18").Value = "3.120.14"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").Value = "'14.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("E21").Value = "  -4.15%  "
$ws.Range("D22").Value = "'440.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("E23").Value = "  +6.39%  "
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("E25").Value = "  -3.14%  "
$ws.Range("D26").Value = "'88.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.27%  "
$ws.Range("D27").Value = "'12.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "'9.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.69%  "
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("D32").Value = "'0.206"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +18.91%  "
$ws.Range("D33").Value = "'26.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.29%  "
$ws.Range("D34").Value = "'0.895"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -10.36%  "
$ws.Range("E35").Value = "  +4.47%  "
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("D37").Value = "'508.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.98%  "
$ws.Range("D38").Value = "'1.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("E40").Value = "  -2.10%  "
$ws.Range("D41").Value = "'0.0893"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.08%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").Value = "'22.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D45").Value = "'3.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +50.73%  "
$ws.Range("E46").Value = "  -2.42%  "
$ws.Range("D47").Value = "'151.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("E48").Value = "  +5.80%  "
$ws.Range("D49").Value = "'45.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("E50").Value = "  -0.51%  "
$ws.Range("E51").Value = "  +1.41%  "
